$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 73.71429000000001
$ws.Range("I38").Value = 74.818184
$ws.Range("J38").Value = 69.666664
$ws.Range("K38").Value = 224.454552
$ws.Range("L38").Value = 208.999992
$ws.Range("M38").Value = 147.545448
$ws.Range("N38").Value = -952.999992

$ws.Range("H43").Value = 1125
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1166.6666
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1166.6666
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -1304.6666

$ws.Range("H129").Value = 688.4286
$ws.Range("I129").Value = 531.3333
$ws.Range("J129").Value = 806.25
$ws.Range("K129").Value = 1593.9999
$ws.Range("L129").Value = 2418.75
$ws.Range("M129").Value = 3406.0001
$ws.Range("N129").Value = -12418.75

$ws.Range("H137").Value = 16951260
$ws.Range("I137").Value = 1110.1
$ws.Range("J137").Value = 52635784
$ws.Range("K137").Value = 3330.3
$ws.Range("L137").Value = 157907352
$ws.Range("M137").Value = -780.2999999999997
$ws.Range("N137").Value = -157912452

$ws.Range("H138").Value = 3266.5398
$ws.Range("I138").Value = 2715.2104
$ws.Range("J138").Value = 4104.56
$ws.Range("K138").Value = 8145.6312
$ws.Range("L138").Value = 12313.68
$ws.Range("M138").Value = -3005.6312
$ws.Range("N138").Value = -22593.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 250
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -138
$ws.Range("N5").Value = -474

$ws.Range("H32").Value = 15112.526
$ws.Range("I32").Value = 12642.676
$ws.Range("J32").Value = 106497
$ws.Range("K32").Value = 12642.676
$ws.Range("L32").Value = 106497
$ws.Range("M32").Value = -12355.676
$ws.Range("N32").Value = -107071

$ws.Range("H37").Value = 9400.888999999999
$ws.Range("I37").Value = 2034
$ws.Range("J37").Value = 10321.75
$ws.Range("K37").Value = 2034
$ws.Range("L37").Value = 10321.75
$ws.Range("M37").Value = -1761
$ws.Range("N37").Value = -10867.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -135
$ws.Range("N4").Value = -480

$ws.Range("H134").Value = 50004028
$ws.Range("I134").Value = 50004028
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 150012084
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -150009549
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 142000.8
$ws.Range("I6").Value = 52501
$ws.Range("J6").Value = 500000
$ws.Range("K6").Value = 52501
$ws.Range("L6").Value = 500000
$ws.Range("M6").Value = -52388
$ws.Range("N6").Value = -500226

$ws.Range("H7").Value = 48.875
$ws.Range("I7").Value = 45
$ws.Range("J7").Value = 52.75
$ws.Range("K7").Value = 45
$ws.Range("L7").Value = 52.75
$ws.Range("M7").Value = 68
$ws.Range("N7").Value = -278.75

$ws.Range("H31").Value = 4388.686
$ws.Range("I31").Value = 1231.8462
$ws.Range("J31").Value = 7671.8
$ws.Range("K31").Value = 1231.8462
$ws.Range("L31").Value = 7671.8
$ws.Range("M31").Value = -936.8462
$ws.Range("N31").Value = -8261.799999999999

$ws.Range("H34").Value = 4388.686
$ws.Range("I34").Value = 1231.8462
$ws.Range("J34").Value = 7671.8
$ws.Range("K34").Value = 1231.8462
$ws.Range("L34").Value = 7671.8
$ws.Range("M34").Value = -1029.8462
$ws.Range("N34").Value = -8075.8

$ws.Range("H41").Value = 7700
$ws.Range("I41").Value = 2550
$ws.Range("J41").Value = 18000
$ws.Range("K41").Value = 2550
$ws.Range("L41").Value = 18000
$ws.Range("M41").Value = -2122
$ws.Range("N41").Value = -18856

$ws.Range("H50").Value = 10836.8
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 10836.8
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 10836.8
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -12086.8

$ws.Range("H51").Value = 9831
$ws.Range("I51").Value = 8800
$ws.Range("J51").Value = 10125.571
$ws.Range("K51").Value = 8800
$ws.Range("L51").Value = 10125.571
$ws.Range("M51").Value = -8064
$ws.Range("N51").Value = -11597.571

$ws.Range("H59").Value = 16056.5
$ws.Range("I59").Value = 16000
$ws.Range("J59").Value = 16064.571
$ws.Range("K59").Value = 16000
$ws.Range("L59").Value = 16064.571
$ws.Range("M59").Value = -14855
$ws.Range("N59").Value = -18354.571

$ws.Range("H60").Value = 10201.714
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 10201.714
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").Value = 10201.714
$ws.Range("N60").Value = -11223.714

$ws.Range("H61").Value = 9831
$ws.Range("I61").Value = 8800
$ws.Range("J61").Value = 10125.571
$ws.Range("K61").Value = 8800
$ws.Range("L61").Value = 10125.571
$ws.Range("M61").Value = -8452
$ws.Range("N61").Value = -10821.571

$ws.Range("H68").Value = 18458.5
$ws.Range("I68").Value = 14268
$ws.Range("J68").Value = 19296.6
$ws.Range("K68").Value = 14268
$ws.Range("L68").Value = 19296.6
$ws.Range("M68").Value = -13519
$ws.Range("N68").Value = -20794.6

$ws.Range("H71").Value = 18458.5
$ws.Range("I71").Value = 14268
$ws.Range("J71").Value = 19296.6
$ws.Range("K71").Value = 42804
$ws.Range("L71").Value = 57889.8
$ws.Range("M71").Value = -39060
$ws.Range("N71").Value = -65377.8

$ws.Range("H74").Value = 14399.777
$ws.Range("I74").Value = 5185
$ws.Range("J74").Value = 17032.572
$ws.Range("K74").Value = 5185
$ws.Range("L74").Value = 17032.572
$ws.Range("M74").Value = -4311
$ws.Range("N74").Value = -18780.572

$ws.Range("H77").Value = 14399.777
$ws.Range("I77").Value = 5185
$ws.Range("J77").Value = 17032.572
$ws.Range("K77").Value = 15555
$ws.Range("L77").Value = 51097.716
$ws.Range("M77").Value = -11187
$ws.Range("N77").Value = -59833.716

$ws.Range("H93").Value = 11664.909
$ws.Range("I93").Value = 9831.4
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 9831.4
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -7959.4
$ws.Range("N93").Value = -33744

$ws.Range("H122").Value = 10223.818
$ws.Range("I122").Value = 11046.2
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 33138.60000000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -30688.60000000001
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 3099.8572
$ws.Range("I132").Value = 2964.5454
$ws.Range("J132").Value = 3596
$ws.Range("K132").Value = 8893.636200000001
$ws.Range("L132").Value = 10788
$ws.Range("M132").Value = -6363.636200000001
$ws.Range("N132").Value = -15848

$ws.Range("H134").Value = 2113.4468
$ws.Range("I134").Value = 1830.9166
$ws.Range("J134").Value = 3038.0908
$ws.Range("K134").Value = 5492.7498
$ws.Range("L134").Value = 9114.2724
$ws.Range("M134").Value = -2957.7498
$ws.Range("N134").Value = -14184.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 39688484
$ws.Range("I131").Value = 111131130
$ws.Range("J131").Value = 27781378
$ws.Range("K131").Value = 333393390
$ws.Range("L131").Value = 83344134
$ws.Range("M131").Value = -333388350
$ws.Range("N131").Value = -83354214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1901.125
$ws.Range("I97").Value = 1744.1428
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1744.1428
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1248.1428
$ws.Range("N97").Value = -3992

$ws.Range("H113").Value = 3081.8333
$ws.Range("I113").Value = 3397.75
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 3397.75
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = -1227.75
$ws.Range("N113").Value = -6790

$ws.Range("H132").Value = 2415.7144
$ws.Range("I132").Value = 1382.4
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4147.200000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1617.200000000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1012.2069
$ws.Range("I46").Value = 1003.0909
$ws.Range("J46").Value = 1040.8572
$ws.Range("K46").Value = 1003.0909
$ws.Range("L46").Value = 1040.8572
$ws.Range("M46").Value = -815.0909
$ws.Range("N46").Value = -1416.8572

$ws.Range("H87").Value = 62189
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 62189
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 62189
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -64435

$ws.Range("H90").Value = 62189
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 62189
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 186567
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -197799

$ws.Range("H93").Value = 1169.5483
$ws.Range("I93").Value = 1287.5264
$ws.Range("J93").Value = 982.75
$ws.Range("K93").Value = 1287.5264
$ws.Range("L93").Value = 982.75
$ws.Range("M93").Value = -39.52639999999997
$ws.Range("N93").Value = -3478.75

$ws.Range("H100").Value = 2340
$ws.Range("I100").Value = 2400
$ws.Range("J100").Value = 2316
$ws.Range("K100").Value = 2400
$ws.Range("L100").Value = 2316
$ws.Range("M100").Value = -1859
$ws.Range("N100").Value = -3398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 375.1
$ws.Range("I100").Value = 364.42856
$ws.Range("J100").Value = 400
$ws.Range("K100").Value = 728.85712
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = -187.85712
$ws.Range("N100").Value = -1882

$ws.Range("H107").Value = 637.2857
$ws.Range("I107").Value = 451.94116
$ws.Range("J107").Value = 1425
$ws.Range("K107").Value = 1355.82348
$ws.Range("L107").Value = 4275
$ws.Range("M107").Value = 564.17652
$ws.Range("N107").Value = -8115
